$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")
$ws.Range("A1:K6").ClearContents() | Out-Null

$labels = @("[0,5)","[5,15)","[15,25)","[25,35)","[35,45)","[45,55)","[55,65)","[65,75)","[75,85)","[85,95)","[95,100]")
$bVals = @(1,0,8,28,56,70,56,28,8,0,1)
$cVals = @(2,0,14,42,70,0,70,42,14,0,2)
$dVals = @(2,0,16,50,60,0,60,50,16,0,2)
$eVals = @(5,0,38,64,21,0,21,64,38,0,5)

for ($i = 0; $i -le 10; $i++) {
    $row = 1 + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}

$jVals = @(5,38,64,21,0,21,64,38,5)
$kVals = @(1,0.76994654383738415,0.66666666666666663,0.58043062325516626,0.5,0.4195693767448338,0.33333333333333337,0.23005345616261588,0)
for ($i = 0; $i -le 8; $i++) {
    $row = 1 + $i
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
    $ws.Cells.Item($row, 11).Value = $kVals[$i]
}

# Second little table (rows 17-21)
$labels2 = @("[50,60)","[60,70)","[70,80)","[80,90)","[90,100)")
for ($i = 0; $i -le 4; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 1).Value = $labels2[$i]
}

$ws.Range("B17").Formula = "=SUM(B5:B7)"
$ws.Range("C17").Formula = "=SUM(C5:C7)"
$ws.Range("D17").Formula = "=SUM(D5:D7)"
$ws.Range("E17").Formula = "=SUM(E5:E7)"

$ws.Range("B18").Formula = "=B8*2"
$ws.Range("C18").Formula = "=C8*2"
$ws.Range("D18").Formula = "=D8*2"
$ws.Range("E18").Formula = "=E8*2"

$ws.Range("B19").Formula = "=B9*2"
$ws.Range("C19").Formula = "=C9*2"
$ws.Range("D19").Formula = "=D9*2"
$ws.Range("E19").Formula = "=E9*2"

$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0

$ws.Range("B21").Formula = "=B11*2"
$ws.Range("C21").Formula = "=C11*2"
$ws.Range("D21").Formula = "=D11*2"
$ws.Range("E21").Formula = "=E11*2"

Write-Host "done all"
